{"js": "// Each entry is [oldText, newText]; applied in document order.\n// The diff only rewrites the text inside existing <w:t> runs (date title\n// plus 25 \"a\u00d7b=c\" table-cell answers); nothing else in the document moves.\nconst replacements = [\n  [\"2025-09-23 Tuesday\", \"2025-09-24 Wednesday\"],\n  [\"870\u00d73=2610\", \"652\u00d77=4564\"],\n  [\"793\u00d72=1586\", \"776\u00d79=6984\"],\n  [\"893\u00d79=8037\", \"419\u00d79=3771\"],\n  [\"567\u00d73=1701\", \"633\u00d76=3798\"],\n  [\"196\u00d75=980\", \"186\u00d77=1302\"],\n  [\"471\u00d79=4239\", \"377\u00d78=3016\"],\n  [\"299\u00d79=2691\", \"884\u00d79=7956\"],\n  [\"479\u00d75=2395\", \"864\u00d73=2592\"],\n  [\"390\u00d79=3510\", \"571\u00d72=1142\"],\n  [\"760\u00d77=5320\", \"273\u00d72=546\"],\n  [\"912\u00d74=3648\", \"739\u00d73=2217\"],\n  [\"576\u00d75=2880\", \"964\u00d74=3856\"],\n  [\"547\u00d73=1641\", \"286\u00d76=1716\"],\n  [\"334\u00d73=1002\", \"960\u00d78=7680\"],\n  [\"973\u00d74=3892\", \"891\u00d74=3564\"],\n  [\"205\u00d77=1435\", \"899\u00d74=3596\"],\n  [\"320\u00d72=640\", \"104\u00d73=312\"],\n  [\"434\u00d72=868\", \"146\u00d73=438\"],\n  [\"403\u00d76=2418\", \"137\u00d78=1096\"],\n  [\"894\u00d78=7152\", \"943\u00d76=5658\"],\n  [\"700\u00d73=2100\", \"554\u00d79=4986\"],\n  [\"957\u00d77=6699\", \"281\u00d75=1405\"],\n  [\"147\u00d77=1029\", \"473\u00d74=1892\"],\n  [\"760\u00d79=6840\", \"259\u00d74=1036\"],\n  [\"647\u00d75=3235\", \"423\u00d79=3807\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replace just this run's text, preserving its formatting (rFonts/sz/etc.).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each entry is old text -> new text, applied in document order.\n# The diff only rewrites the text inside existing runs (date title plus\n# 25 \"a\u00d7b=c\" table-cell answers); nothing else in the document moves.\n$replacements = @(\n    @('2025-09-23 Tuesday', '2025-09-24 Wednesday'),\n    @('870\u00d73=2610', '652\u00d77=4564'),\n    @('793\u00d72=1586', '776\u00d79=6984'),\n    @('893\u00d79=8037', '419\u00d79=3771'),\n    @('567\u00d73=1701', '633\u00d76=3798'),\n    @('196\u00d75=980', '186\u00d77=1302'),\n    @('471\u00d79=4239', '377\u00d78=3016'),\n    @('299\u00d79=2691', '884\u00d79=7956'),\n    @('479\u00d75=2395', '864\u00d73=2592'),\n    @('390\u00d79=3510', '571\u00d72=1142'),\n    @('760\u00d77=5320', '273\u00d72=546'),\n    @('912\u00d74=3648', '739\u00d73=2217'),\n    @('576\u00d75=2880', '964\u00d74=3856'),\n    @('547\u00d73=1641', '286\u00d76=1716'),\n    @('334\u00d73=1002', '960\u00d78=7680'),\n    @('973\u00d74=3892', '891\u00d74=3564'),\n    @('205\u00d77=1435', '899\u00d74=3596'),\n    @('320\u00d72=640', '104\u00d73=312'),\n    @('434\u00d72=868', '146\u00d73=438'),\n    @('403\u00d76=2418', '137\u00d78=1096'),\n    @('894\u00d78=7152', '943\u00d76=5658'),\n    @('700\u00d73=2100', '554\u00d79=4986'),\n    @('957\u00d77=6699', '281\u00d75=1405'),\n    @('147\u00d77=1029', '473\u00d74=1892'),\n    @('760\u00d79=6840', '259\u00d74=1036'),\n    @('647\u00d75=3235', '423\u00d79=3807'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Search the whole document body and replace the (unique) matching run's\n    # text in place, which preserves its existing character formatting.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Replace = 2 -> wdReplaceAll\n    $result = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
